$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Running C and D suite only: A Suite -> N, C Suite -> Y, D Suite -> Y
$ws.Range("C2").Value = "N"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"

# Update active selection/cell as seen in diff
$ws.Range("C5").Select()
